$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)
$ws.Columns("E:E").Select()
$ws.Columns("E:E").Delete()
$lo.Resize($ws.Range("A1:K92"))
